$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.085.81'
$ws.Range('E2').Value = '  -1.99%  '
$ws.Range('D3').Value = '1.834.88'
$ws.Range('E3').Value = '  -0.89%  '
$ws.Range('D4').Value = '''1.000'
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').Value = '''324.24'
$ws.Range('E5').Value = '  -2.95%  '
$ws.Range('D6').Value = '''0.9996'
$ws.Range('E6').Value = '  -0.15%  '
$ws.Range('D7').Value = '''0.4644'
$ws.Range('E7').Value = '  -0.25%  '
$ws.Range('D8').Value = '''0.3870'
$ws.Range('E8').Value = '  -1.05%  '
$ws.Range('D9').Value = '''0.07873'
$ws.Range('E9').Value = '  -0.66%  '
$ws.Range('D10').Value = '''0.9611'
$ws.Range('E10').Value = '  -2.40%  '
$ws.Range('D11').Value = '''21.94'
$ws.Range('E11').Value = '  -1.64%  '
$ws.Range('D12').Value = '1.852.51'
$ws.Range('E12').Value = '  -4.06%  '
$ws.Range('D13').Value = '''5.690'
$ws.Range('E13').Value = '  -2.70%  '
$ws.Range('D14').Value = '''6.909'
$ws.Range('E14').Value = '  -1.26%  '
$ws.Range('D15').Value = '''0.06855'
$ws.Range('E15').Value = '  +0.00%  '
$ws.Range('D16').Value = '''87.18'
$ws.Range('E16').Value = '  -0.69%  '
$ws.Range('E17').Value = '  -0.15%  '
$ws.Range('D18').Value = '''0.000009946'
$ws.Range('E18').Value = '  -1.45%  '
$ws.Range('D19').Value = '''16.67'
$ws.Range('E19').Value = '  -2.48%  '
$ws.Range('E20').Value = '  +0.02%  '
$ws.Range('D21').Value = '28.090.86'
$ws.Range('E21').Value = '  -2.07%  '
$ws.Range('D22').Value = '''5.327'
$ws.Range('E22').Value = '  -1.28%  '
$ws.Range('D23').Value = '''11.00'
$ws.Range('D24').Value = '''2.096'
$ws.Range('E24').Value = '  -1.48%  '
$ws.Range('D25').Value = '2.034.18'
$ws.Range('E25').Value = '  -5.82%  '
$ws.Range('D26').Value = '''154.03'
$ws.Range('E26').Value = '  +0.66%  '
$ws.Range('D27').Value = '''19.18'
$ws.Range('E27').Value = '  -1.27%  '
$ws.Range('D28').Value = '''5.695'
$ws.Range('E28').Value = '  -6.29%  '
$ws.Range('D29').Value = '''1.968'
$ws.Range('E29').Value = '  -2.76%  '
$ws.Range('D30').Value = '''117.97'
$ws.Range('E30').Value = '  +0.24%  '
$ws.Range('B31').Value = 'ImmutableX'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D31').Value = '''0.9364'
$ws.Range('E31').Value = '  -4.12%  '
$ws.Range('B32').Value = 'Stellar'
$ws.Range('C32').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D32').Value = '''0.09254'
$ws.Range('E32').Value = '  -1.76%  '
$ws.Range('D33').Value = '''5.275'
$ws.Range('E33').Value = '  -1.71%  '
$ws.Range('D34').Value = '''1.323'
$ws.Range('E34').Value = '  -2.05%  '
$ws.Range('D35').Value = '''3.293'
$ws.Range('E35').Value = '  -5.39%  '
$ws.Range('D36').Value = '''0.05850'
$ws.Range('E36').Value = '  -5.34%  '
$ws.Range('D37').Value = '''0.02125'
$ws.Range('E37').Value = '  -3.63%  '
$ws.Range('D38').Value = '''1.143'
$ws.Range('E38').Value = '  -1.80%  '
$ws.Range('D39').Value = '''7.781'
$ws.Range('E39').Value = '  +2.08%  '
$ws.Range('D40').Value = '''0.5601'
$ws.Range('E40').Value = '  -2.14%  '
$ws.Range('D41').Value = '''9.898'
$ws.Range('E41').Value = '  -2.51%  '
$ws.Range('D42').Value = '''0.1765'
$ws.Range('E42').Value = '  -1.93%  '
$ws.Range('D43').Value = '''0.07214'
$ws.Range('E43').Value = '  +0.79%  '
$ws.Range('D44').Value = '''11.63'
$ws.Range('E44').Value = '  -1.23%  '
$ws.Range('D45').Value = '''0.5272'
$ws.Range('E45').Value = '  -2.39%  '
$ws.Range('D46').Value = '''2.119'
$ws.Range('E46').Value = '  -10.95%  '
$ws.Range('D47').Value = '''1.121'
$ws.Range('E47').Value = '  -10.17%  '
$ws.Range('D48').Value = '''1.835'
$ws.Range('E48').Value = '  -3.90%  '
$ws.Range('D49').Value = '''112.79'
$ws.Range('E49').Value = '  -1.09%  '
$ws.Range('D50').Value = '''0.9994'
$ws.Range('E50').Value = '  -0.13%  '
$ws.Range('D51').Value = '''2.320'
$ws.Range('E51').Value = '  +0.04%  '
